$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$updates = @(
  @("C2", -1),
  @("D2", 5),
  @("E2", 0),
  @("C3", -1),
  @("D3", 25),
  @("E3", 0),
  @("I3", 16),
  @("C4", -1),
  @("D4", 3),
  @("E4", 0),
  @("C5", -1),
  @("D5", 9),
  @("E5", 0),
  @("I5", 6),
  @("C6", -1),
  @("D6", 11),
  @("E6", 0),
  @("I6", 10),
  @("C7", -1),
  @("D7", 3),
  @("E7", 0),
  @("I7", 28),
  @("C8", -1),
  @("D8", 21),
  @("E8", 0),
  @("I8", 13),
  @("C9", -1),
  @("D9", 21),
  @("E9", 0),
  @("I9", 13),
  @("C10", -1),
  @("D10", 14),
  @("E10", 0),
  @("I10", 11),
  @("C11", -1),
  @("D11", 5),
  @("E11", 0),
  @("I11", 3),
  @("C12", -1),
  @("D12", 9),
  @("E12", 0),
  @("I12", 5),
  @("C13", -1),
  @("D13", 9),
  @("E13", 0),
  @("I13", 6),
  @("C14", -1),
  @("E14", 0),
  @("C15", -1),
  @("D15", 8),
  @("E15", 0),
  @("I15", 4),
  @("C16", -1),
  @("D16", 3),
  @("E16", 0),
  @("C17", -1),
  @("E17", 0),
  @("C18", -1),
  @("E18", 0),
  @("C19", -1),
  @("D19", 8),
  @("E19", 0),
  @("C20", -1),
  @("E20", 0),
  @("C21", -1),
  @("D21", 13),
  @("E21", 0),
  @("C22", -1),
  @("D22", 18),
  @("E22", 0),
  @("I22", 2),
  @("C23", -1),
  @("D23", 3),
  @("E23", 0),
  @("C24", -1),
  @("D24", 4),
  @("E24", 0),
  @("C25", -1),
  @("D25", 21),
  @("E25", 0),
  @("I25", 13),
  @("C26", -1),
  @("D26", 6),
  @("E26", 0),
  @("C27", -1),
  @("D27", 21),
  @("E27", 0),
  @("I27", 13),
  @("C28", -1),
  @("D28", 14),
  @("E28", 0),
  @("I28", 10),
  @("C29", -1),
  @("D29", 12),
  @("E29", 0),
  @("C30", -1),
  @("D30", 5),
  @("E30", 0),
  @("C31", -1),
  @("D31", 9),
  @("E31", 0),
  @("I31", 3),
  @("C32", -1),
  @("D32", 14),
  @("E32", 0),
  @("I32", 7),
  @("C33", -1),
  @("D33", 7),
  @("E33", 0),
  @("I33", 5),
  @("C34", -1),
  @("D34", 18),
  @("E34", 0),
  @("I34", 13)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("Ochiai")
$updates = @(
  @("C2", -1),
  @("D2", 5),
  @("E2", 0),
  @("C3", -1),
  @("D3", 2),
  @("E3", 0),
  @("C4", -1),
  @("D4", 3),
  @("E4", 0),
  @("C5", -1),
  @("D5", 6),
  @("E5", 0),
  @("I5", 4),
  @("C6", -1),
  @("D6", 6),
  @("E6", 0),
  @("C7", -1),
  @("D7", 3),
  @("E7", 0),
  @("C8", -1),
  @("D8", 5),
  @("E8", 0),
  @("I8", 3),
  @("C9", -1),
  @("D9", 5),
  @("E9", 0),
  @("I9", 3),
  @("C10", -1),
  @("D10", 2),
  @("E10", 0),
  @("C11", -1),
  @("D11", 5),
  @("E11", 0),
  @("I11", 3),
  @("C12", -1),
  @("D12", 9),
  @("E12", 0),
  @("I12", 5),
  @("C13", -1),
  @("D13", 6),
  @("E13", 0),
  @("I13", 4),
  @("C14", -1),
  @("E14", 0),
  @("C15", -1),
  @("D15", 6),
  @("E15", 0),
  @("I15", 4),
  @("C16", -1),
  @("D16", 3),
  @("E16", 0),
  @("C17", -1),
  @("E17", 0),
  @("C18", -1),
  @("E18", 0),
  @("C19", -1),
  @("E19", 0),
  @("C20", -1),
  @("E20", 0),
  @("C21", -1),
  @("D21", 2),
  @("E21", 0),
  @("C22", -1),
  @("D22", 27),
  @("E22", 0),
  @("I22", 14),
  @("C23", -1),
  @("D23", 3),
  @("E23", 0),
  @("C24", -1),
  @("D24", 4),
  @("E24", 0),
  @("C25", -1),
  @("D25", 5),
  @("E25", 0),
  @("I25", 3),
  @("C26", -1),
  @("E26", 0),
  @("C27", -1),
  @("D27", 5),
  @("E27", 0),
  @("I27", 3),
  @("C28", -1),
  @("D28", 13),
  @("E28", 0),
  @("I28", 9),
  @("C29", -1),
  @("D29", 2),
  @("E29", 0),
  @("C30", -1),
  @("E30", 0),
  @("C31", -1),
  @("E31", 0),
  @("C32", -1),
  @("D32", 6),
  @("E32", 0),
  @("I32", 4),
  @("C33", -1),
  @("D33", 7),
  @("E33", 0),
  @("I33", 5),
  @("C34", -1),
  @("D34", 5),
  @("E34", 0),
  @("I34", 3)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("Op2")
$updates = @(
  @("C2", -1),
  @("D2", 5),
  @("E2", 0),
  @("C3", -1),
  @("D3", 2),
  @("E3", 0),
  @("C4", -1),
  @("D4", 3),
  @("E4", 0),
  @("C5", -1),
  @("D5", 6),
  @("E5", 0),
  @("I5", 4),
  @("C6", -1),
  @("D6", 6),
  @("E6", 0),
  @("C7", -1),
  @("D7", 3),
  @("E7", 0),
  @("C8", -1),
  @("D8", 5),
  @("E8", 0),
  @("I8", 3),
  @("C9", -1),
  @("D9", 5),
  @("E9", 0),
  @("I9", 3),
  @("C10", -1),
  @("D10", 2),
  @("E10", 0),
  @("C11", -1),
  @("D11", 5),
  @("E11", 0),
  @("I11", 3),
  @("C12", -1),
  @("D12", 7),
  @("E12", 0),
  @("I12", 5),
  @("C13", -1),
  @("D13", 6),
  @("E13", 0),
  @("I13", 4),
  @("C14", -1),
  @("E14", 0),
  @("C15", -1),
  @("D15", 6),
  @("E15", 0),
  @("I15", 4),
  @("C16", -1),
  @("D16", 3),
  @("E16", 0),
  @("C17", -1),
  @("E17", 0),
  @("C18", -1),
  @("E18", 0),
  @("C19", -1),
  @("E19", 0),
  @("C20", -1),
  @("E20", 0),
  @("C21", -1),
  @("D21", 2),
  @("E21", 0),
  @("C22", -1),
  @("D22", 32),
  @("E22", 0),
  @("I22", 28),
  @("C23", -1),
  @("D23", 3),
  @("E23", 0),
  @("C24", -1),
  @("D24", 4),
  @("E24", 0),
  @("C25", -1),
  @("D25", 5),
  @("E25", 0),
  @("I25", 3),
  @("C26", -1),
  @("E26", 0),
  @("C27", -1),
  @("D27", 5),
  @("E27", 0),
  @("I27", 3),
  @("C28", -1),
  @("D28", 13),
  @("E28", 0),
  @("I28", 9),
  @("C29", -1),
  @("D29", 2),
  @("E29", 0),
  @("C30", -1),
  @("E30", 0),
  @("C31", -1),
  @("E31", 0),
  @("C32", -1),
  @("D32", 6),
  @("E32", 0),
  @("I32", 4),
  @("C33", -1),
  @("D33", 7),
  @("E33", 0),
  @("I33", 5),
  @("C34", -1),
  @("D34", 5),
  @("E34", 0),
  @("I34", 3)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("Barinel")
$updates = @(
  @("C2", -1),
  @("D2", 5),
  @("E2", 0),
  @("C3", -1),
  @("D3", 25),
  @("E3", 0),
  @("I3", 16),
  @("C4", -1),
  @("D4", 3),
  @("E4", 0),
  @("C5", -1),
  @("D5", 9),
  @("E5", 0),
  @("I5", 6),
  @("C6", -1),
  @("D6", 11),
  @("E6", 0),
  @("I6", 10),
  @("C7", -1),
  @("D7", 3),
  @("E7", 0),
  @("I7", 28),
  @("C8", -1),
  @("D8", 21),
  @("E8", 0),
  @("I8", 13),
  @("C9", -1),
  @("D9", 21),
  @("E9", 0),
  @("I9", 13),
  @("C10", -1),
  @("D10", 14),
  @("E10", 0),
  @("I10", 11),
  @("C11", -1),
  @("D11", 5),
  @("E11", 0),
  @("I11", 3),
  @("C12", -1),
  @("D12", 9),
  @("E12", 0),
  @("I12", 5),
  @("C13", -1),
  @("D13", 9),
  @("E13", 0),
  @("I13", 6),
  @("C14", -1),
  @("E14", 0),
  @("C15", -1),
  @("D15", 8),
  @("E15", 0),
  @("I15", 4),
  @("C16", -1),
  @("D16", 3),
  @("E16", 0),
  @("C17", -1),
  @("E17", 0),
  @("C18", -1),
  @("E18", 0),
  @("C19", -1),
  @("D19", 8),
  @("E19", 0),
  @("C20", -1),
  @("E20", 0),
  @("C21", -1),
  @("D21", 13),
  @("E21", 0),
  @("C22", -1),
  @("D22", 18),
  @("E22", 0),
  @("I22", 2),
  @("C23", -1),
  @("D23", 3),
  @("E23", 0),
  @("C24", -1),
  @("D24", 4),
  @("E24", 0),
  @("C25", -1),
  @("D25", 21),
  @("E25", 0),
  @("I25", 13),
  @("C26", -1),
  @("D26", 6),
  @("E26", 0),
  @("C27", -1),
  @("D27", 21),
  @("E27", 0),
  @("I27", 13),
  @("C28", -1),
  @("D28", 14),
  @("E28", 0),
  @("I28", 10),
  @("C29", -1),
  @("D29", 12),
  @("E29", 0),
  @("C30", -1),
  @("D30", 5),
  @("E30", 0),
  @("C31", -1),
  @("D31", 9),
  @("E31", 0),
  @("I31", 3),
  @("C32", -1),
  @("D32", 14),
  @("E32", 0),
  @("I32", 7),
  @("C33", -1),
  @("D33", 7),
  @("E33", 0),
  @("I33", 5),
  @("C34", -1),
  @("D34", 18),
  @("E34", 0),
  @("I34", 13)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("Dstar")
$updates = @(
  @("C2", -1),
  @("D2", 5),
  @("E2", 0),
  @("C3", -1),
  @("D3", 2),
  @("E3", 0),
  @("C4", -1),
  @("D4", 3),
  @("E4", 0),
  @("C5", -1),
  @("D5", 6),
  @("E5", 0),
  @("I5", 4),
  @("C6", -1),
  @("D6", 6),
  @("E6", 0),
  @("C7", -1),
  @("D7", 3),
  @("E7", 0),
  @("C8", -1),
  @("D8", 5),
  @("E8", 0),
  @("I8", 3),
  @("C9", -1),
  @("D9", 5),
  @("E9", 0),
  @("I9", 3),
  @("C10", -1),
  @("D10", 2),
  @("E10", 0),
  @("C11", -1),
  @("D11", 5),
  @("E11", 0),
  @("I11", 3),
  @("C12", -1),
  @("D12", 7),
  @("E12", 0),
  @("I12", 5),
  @("C13", -1),
  @("D13", 6),
  @("E13", 0),
  @("I13", 4),
  @("C14", -1),
  @("E14", 0),
  @("C15", -1),
  @("D15", 6),
  @("E15", 0),
  @("I15", 4),
  @("C16", -1),
  @("D16", 3),
  @("E16", 0),
  @("C17", -1),
  @("E17", 0),
  @("C18", -1),
  @("E18", 0),
  @("C19", -1),
  @("E19", 0),
  @("C20", -1),
  @("E20", 0),
  @("C21", -1),
  @("D21", 2),
  @("E21", 0),
  @("C22", -1),
  @("D22", 32),
  @("E22", 0),
  @("I22", 14),
  @("C23", -1),
  @("D23", 3),
  @("E23", 0),
  @("C24", -1),
  @("D24", 4),
  @("E24", 0),
  @("C25", -1),
  @("D25", 5),
  @("E25", 0),
  @("I25", 3),
  @("C26", -1),
  @("E26", 0),
  @("C27", -1),
  @("D27", 5),
  @("E27", 0),
  @("I27", 3),
  @("C28", -1),
  @("D28", 13),
  @("E28", 0),
  @("I28", 9),
  @("C29", -1),
  @("D29", 2),
  @("E29", 0),
  @("C30", -1),
  @("E30", 0),
  @("C31", -1),
  @("E31", 0),
  @("C32", -1),
  @("D32", 6),
  @("E32", 0),
  @("I32", 4),
  @("C33", -1),
  @("D33", 7),
  @("E33", 0),
  @("I33", 5),
  @("C34", -1),
  @("D34", 5),
  @("E34", 0),
  @("I34", 3)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }
